# Last Sync: 2021-04-06 10:30:02
#
# Rewrites the "Validation" sheet's two small lookup tables (rows 1-3 and
# rows 7-9) to use the new, shortened category labels (e.g. "SOW" instead
# of "SOW - Client" / "SOW - Vendor") and widens each table with the
# additional document-type columns that were introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header (B1:O1) -------------------------------------------------
$ws.Range("B1").Value = "SOW"
$ws.Range("C1").Value = "MSA"
$ws.Range("D1").Value = "Hubspot legal terms & conditions"
$ws.Range("E1").Value = "NDA"
$ws.Range("F1").Value = "Code of Conduct"
$ws.Range("G1").Value = "BAA"
$ws.Range("H1").Value = "PSA"
$ws.Range("I1").Value = "Miscellaneous"
$ws.Range("J1").Value = "Order Form"
$ws.Range("K1").Value = "DPA"
$ws.Range("L1").Value = "SDPA"
$ws.Range("M1").Value = "Partner Affiliate Adoption Agreement"
$ws.Range("N1").Value = "Referral Agreement"
$ws.Range("O1").Value = "Engagement Letter"

# New header cells need the same bold/bordered/centered style as the rest
# of the header row (copy formats from the already-styled L1 cell).
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1:O1").PasteSpecial(-4122) | Out-Null

# --- Row 2 (Q-Tempelate) ---------------------------------------------------
$ws.Range("A2").Value = "Q-Tempelate"
$ws.Range("B2").Value = 1.125
$ws.Range("C2").Value = 4.466666666666667
$ws.Range("D2").Value = "Not enough data"
$ws.Range("E2").Value = "Not enough data"
$ws.Range("F2").Value = "Not enough data"
$ws.Range("G2").Value = "Not enough data"
$ws.Range("H2").Value = "Not enough data"
$ws.Range("I2").Value = "Not enough data"
$ws.Range("J2").Value = "Not enough data"
$ws.Range("K2").Value = "Not enough data"
$ws.Range("L2").Value = "Not enough data"
$ws.Range("M2").Value = "Not enough data"
$ws.Range("N2").Value = "Not enough data"
$ws.Range("O2").Value = "Not enough data"

# --- Row 3 (Non-Q Tempelate) -----------------------------------------------
$ws.Range("A3").Value = "Non-Q Tempelate"
$ws.Range("B3").Value = 2.5
$ws.Range("C3").Value = 3.125
$ws.Range("D3").Value = "Not enough data"
$ws.Range("E3").Value = 0.6
$ws.Range("F3").Value = "Not enough data"
$ws.Range("G3").Value = "Not enough data"
$ws.Range("H3").Value = 2.333333333333333
$ws.Range("I3").Value = "Not enough data"
$ws.Range("J3").Value = "Not enough data"
$ws.Range("K3").Value = "Not enough data"
$ws.Range("L3").Value = "Not enough data"
$ws.Range("M3").Value = "Not enough data"
$ws.Range("N3").Value = "Not enough data"
$ws.Range("O3").Value = "Not enough data"

# --- Row 7 header (B7:Q7) ---------------------------------------------------
$ws.Range("B7").Value = "SOW"
$ws.Range("C7").Value = "ZoomInfo Recurring credits"
$ws.Range("D7").Value = "Change Order"
$ws.Range("E7").Value = "NDA"
$ws.Range("F7").Value = "PSA"
$ws.Range("G7").Value = "Miscellaneous"
$ws.Range("H7").Value = "Purchase Order"
$ws.Range("I7").Value = "Lease Agreement"
$ws.Range("J7").Value = "software license agreement"
$ws.Range("K7").Value = "POC"
$ws.Range("L7").Value = "MSA"
$ws.Range("M7").Value = "Sponsorship Agreement"
$ws.Range("N7").Value = "Partner Affliate Consent Letter"
$ws.Range("O7").Value = "MSA "
$ws.Range("P7").Value = "Addendum"
$ws.Range("Q7").Value = "Partner Agreement"

# New cells in the second table's header row also need the header style
# (L7 through Q7 did not exist before).
$ws.Range("K7").Copy() | Out-Null
$ws.Range("L7:Q7").PasteSpecial(-4122) | Out-Null

# --- Row 8 (Q-Tempelate) ---------------------------------------------------
$ws.Range("A8").Value = "Q-Tempelate"
$ws.Range("B8").Value = 1.285714285714286
$ws.Range("C8").Value = "Not enough data"
$ws.Range("D8").Value = "Not enough data"
$ws.Range("E8").Value = "Not enough data"
$ws.Range("F8").Value = "Not enough data"
$ws.Range("G8").Value = "Not enough data"
$ws.Range("H8").Value = "Not enough data"
$ws.Range("I8").Value = "Not enough data"
$ws.Range("J8").Value = "Not enough data"
$ws.Range("K8").Value = "Not enough data"
$ws.Range("L8").Value = "Not enough data"
$ws.Range("M8").Value = "Not enough data"
$ws.Range("N8").Value = "Not enough data"
$ws.Range("O8").Value = "Not enough data"
$ws.Range("P8").Value = "Not enough data"
$ws.Range("Q8").Value = "Not enough data"

# --- Row 9 (Non-Q Tempelate) -----------------------------------------------
$ws.Range("A9").Value = "Non-Q Tempelate"
$ws.Range("B9").Value = 1.75
$ws.Range("C9").Value = "Not enough data"
$ws.Range("D9").Value = "Not enough data"
$ws.Range("E9").Value = 1.818181818181818
$ws.Range("F9").Value = "Not enough data"
$ws.Range("G9").Value = "Not enough data"
$ws.Range("H9").Value = "Not enough data"
$ws.Range("I9").Value = "Not enough data"
$ws.Range("J9").Value = "Not enough data"
$ws.Range("K9").Value = "Not enough data"
$ws.Range("L9").Value = "Not enough data"
$ws.Range("M9").Value = "Not enough data"
$ws.Range("N9").Value = "Not enough data"
$ws.Range("O9").Value = "Not enough data"
$ws.Range("P9").Value = "Not enough data"
$ws.Range("Q9").Value = "Not enough data"

Write-Host "Validation sheet rewritten"
